$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (Nombre, Apellido, DNI) for "Pipo Pescador", DNI 666.
$ws.Range("A7").Value = "Pipo"
$ws.Range("B7").Value = "Pescador"
$ws.Range("C7").Value = 666
